$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Purchase 22-23")

# The "6" line (row 18) used to be the last entry of its group, with F18
# totalling the group via "=E18". A new invoice line is being added to that
# same group, so insert a fresh row above the following blank separator
# (row 19), which pushes the separator and the "7" group down by one.
$ws.Rows.Item(19).Insert() | Out-Null

# Pull the row's formatting (number formats, fonts, borders, row height) down
# from row 18 so the new line looks like the rest of the table.
$ws.Range("A18:F18").Copy($ws.Range("A19:F19")) | Out-Null
$ws.Rows.Item(19).RowHeight = $ws.Rows.Item(18).RowHeight

# Fill in the new invoice line.
$ws.Range("A19").Value = $null
$ws.Range("B19").Value = 45348
$ws.Range("C19").Value = 1693
$ws.Range("D19").Value = "Aquachemitech"
$ws.Range("E19").Value = 28000
$ws.Range("F19").Formula = "=E18+E19"

# Row 18 is no longer the last row of its group, so it no longer carries the
# running subtotal - that now lives on row 19.
$ws.Range("F18").Formula = ""

# Mirror the author's final selection.
$ws.Range("F18").Select() | Out-Null
